$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------
# 1) Remove the GLWB and VUL sheets
# -----------------------------------------------------------------
$wb.Worksheets.Item("GLWB").Delete()
$wb.Worksheets.Item("VUL").Delete()

# -----------------------------------------------------------------
# 2) ParamList sheet: drop asmp_file / sensitivity rows, rename two
#    parameter ids, and add asmp_id / mp_file_id rows
# -----------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("ParamList")

# Delete "asmp_file" (row 6) and "sensitivity" (row 11, becomes row 10
# after the first delete) - delete bottom-up so row numbers stay valid.
$ws1.Rows.Item(11).Delete()
$ws1.Rows.Item(6).Delete()

# Rename model_point_file_stem -> model_point_file_prefix (row 4)
$ws1.Range("A4").Value = "model_point_file_prefix"
# Rename int_rate_prefix -> scen_file_prefix (was row 7, now row 6)
$ws1.Range("A6").Value = "scen_file_prefix"

# Insert two new rows before "expense_acq" (currently row 12) for the
# new mp_file_id / asmp_id run parameters
$ws1.Rows.Item(12).Insert()
$ws1.Rows.Item(12).Insert()
$ws1.Range("A12").Value = "mp_file_id"
$ws1.Range("B12").Value = "RUN"
$ws1.Range("C12").Value = "Model point file ID"
$ws1.Range("A13").Value = "asmp_id"
$ws1.Range("B13").Value = "RUN"
$ws1.Range("C13").Value = "Assumption file ID"

$ws1.Range("C13").Select()

# -----------------------------------------------------------------
# 3) ConstParams sheet: drop base_date row, rename two parameter ids
# -----------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("ConstParams")

$ws2.Rows.Item(2).Delete()
# model_point_file_stem -> model_point_file_prefix (was row4, now row3)
$ws2.Range("A3").Value = "model_point_file_prefix"
# int_rate_prefix -> scen_file_prefix (was row8, now row7)
$ws2.Range("A7").Value = "scen_file_prefix"

$ws2.Range("A8").Select()

# -----------------------------------------------------------------
# 4) RunParams sheet: remove scen_file/asmp_file/disc_file/sensitivity
#    columns, reorder date_id/asmp_id/sens_int_rate and add mp_file_id
# -----------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("RunParams")

# Remove columns C:F (scen_file, asmp_file, disc_file, sensitivity)
$ws3.Range("C1:F1").EntireColumn.Delete()
# Now columns are: A run_id, B base_date, C sens_int_rate, D date_id,
# E asmp_id, F description

# Move date_id (D) before sens_int_rate (C)
$ws3.Columns.Item(4).Cut()
$ws3.Columns.Item(3).Insert()
# Now: A run_id, B base_date, C date_id, D sens_int_rate, E asmp_id,
# F description

# Move sens_int_rate (D) after asmp_id (E)
$ws3.Columns.Item(4).Cut()
$ws3.Columns.Item(6).Insert()
# Now: A run_id, B base_date, C date_id, D asmp_id, E sens_int_rate,
# F description

# Insert a new blank column for mp_file_id before asmp_id
$ws3.Columns.Item(4).Insert()
# Now: A run_id, B base_date, C date_id, D (blank), E asmp_id,
# F sens_int_rate, G description

$ws3.Columns.Item(4).ColumnWidth = 15.55
$ws3.Columns.Item(5).ColumnWidth = 10.93

$ws3.Range("D1").Value = "mp_file_id"

$ws3.Range("D2").Value = "202401NB"
$ws3.Range("D3").Value = "2023Q4IF"
$ws3.Range("D4").Value = "2023Q4IF"
$ws3.Range("D5").Value = "2023Q4IF"
$ws3.Range("D6").Value = "2022Q4IF"

$ws3.Range("E28").Select()

# -----------------------------------------------------------------
# 5) Cosmetic: selection changes recorded on the other sheets
# -----------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("SpaceParams")
$ws4.Range("E4").Select()

$ws5 = $wb.Worksheets.Item("GMXB")
$ws5.Range("G26").Select()

$ws1.Activate()
